$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3294.875
$ws.Range("I113").Value = 2850.5881
$ws.Range("J113").Value = 3798.4
$ws.Range("K113").Value = 2850.5881
$ws.Range("L113").Value = 3798.4
$ws.Range("M113").Value = 403.4119000000001
$ws.Range("N113").Value = -10306.4

$ws.Range("H132").Value = 3848206.8
$ws.Range("I132").Value = 1956.5933
$ws.Range("J132").Value = 41669668
$ws.Range("K132").Value = 5869.7799
$ws.Range("L132").Value = 125009004
$ws.Range("M132").Value = -3339.7799
$ws.Range("N132").Value = -125014064

$ws.Range("H138").Value = 4019662.8
$ws.Range("I138").Value = 10418863
$ws.Range("J138").Value = 4478.0195
$ws.Range("K138").Value = 31256589
$ws.Range("L138").Value = 13434.0585
$ws.Range("M138").Value = -31251449
$ws.Range("N138").Value = -23714.0585

$ws.Range("H141").Value = 1801
$ws.Range("I141").Value = 1313.75
$ws.Range("J141").Value = 3750
$ws.Range("K141").Value = 3941.25
$ws.Range("L141").Value = 11250
$ws.Range("M141").Value = 1238.75
$ws.Range("N141").Value = -21610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18967.453
$ws.Range("I32").Value = 18983.617
$ws.Range("J32").Value = 18725
$ws.Range("K32").Value = 18983.617
$ws.Range("L32").Value = 18725
$ws.Range("M32").Value = -18696.617
$ws.Range("N32").Value = -19299

$ws.Range("H37").Value = 8042.8
$ws.Range("J37").Value = 8042.8
$ws.Range("L37").Value = 8042.8
$ws.Range("N37").Value = -8588.799999999999

$ws.Range("H55").Value = 31269.182
$ws.Range("J55").Value = 34096.1
$ws.Range("L55").Value = 34096.1
$ws.Range("N55").Value = -34726.1

$ws.Range("H61").Value = 1646.8064
$ws.Range("I61").Value = 1665.6394
$ws.Range("J61").Value = 498
$ws.Range("K61").Value = 1665.6394
$ws.Range("L61").Value = 498
$ws.Range("M61").Value = -1453.6394
$ws.Range("N61").Value = -922

$ws.Range("H74").Value = 5899.5454
$ws.Range("I74").Value = 1314.5555
$ws.Range("J74").Value = 26532
$ws.Range("K74").Value = 1314.5555
$ws.Range("L74").Value = 26532
$ws.Range("M74").Value = -440.5554999999999
$ws.Range("N74").Value = -28280

$ws.Range("H77").Value = 5899.5454
$ws.Range("I77").Value = 1314.5555
$ws.Range("J77").Value = 26532
$ws.Range("K77").Value = 6572.7775
$ws.Range("L77").Value = 132660
$ws.Range("M77").Value = -2204.7775
$ws.Range("N77").Value = -141396

$ws.Range("H80").Value = 27580
$ws.Range("J80").Value = 28406.666
$ws.Range("L80").Value = 28406.666
$ws.Range("N80").Value = -30402.666

$ws.Range("H83").Value = 27580
$ws.Range("J83").Value = 28406.666
$ws.Range("L83").Value = 85219.99800000001
$ws.Range("N83").Value = -95203.99800000001

$ws.Range("H122").Value = 11457.467
$ws.Range("I122").Value = 15296.4
$ws.Range("J122").Value = 3779.6
$ws.Range("K122").Value = 45889.2
$ws.Range("L122").Value = 11338.8
$ws.Range("M122").Value = -43439.2
$ws.Range("N122").Value = -16238.8

$ws.Range("H132").Value = 2015.2195
$ws.Range("I132").Value = 1506.2059
$ws.Range("K132").Value = 4518.6177
$ws.Range("M132").Value = -1988.6177

$ws.Range("H136").Value = 1646.8064
$ws.Range("I136").Value = 1665.6394
$ws.Range("J136").Value = 498
$ws.Range("K136").Value = 4996.9182
$ws.Range("L136").Value = 1494
$ws.Range("M136").Value = -2446.9182
$ws.Range("N136").Value = -6594

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 677.7083
$ws.Range("I64").Value = 387.25
$ws.Range("J64").Value = 968.1667
$ws.Range("K64").Value = 387.25
$ws.Range("L64").Value = 968.1667
$ws.Range("M64").Value = -162.25
$ws.Range("N64").Value = -1418.1667

$ws.Range("H67").Value = 677.7083
$ws.Range("I67").Value = 387.25
$ws.Range("J67").Value = 968.1667
$ws.Range("K67").Value = 387.25
$ws.Range("L67").Value = 968.1667
$ws.Range("M67").Value = 392.75
$ws.Range("N67").Value = -2528.1667

$ws.Range("H105").Value = 4770.222
$ws.Range("I105").Value = 4174.5
$ws.Range("J105").Value = 4873.826
$ws.Range("K105").Value = 4174.5
$ws.Range("L105").Value = 4873.826
$ws.Range("M105").Value = -2427.5
$ws.Range("N105").Value = -8367.826000000001

$ws.Range("H134").Value = 4227.298
$ws.Range("I134").Value = 3319.4546
$ws.Range("J134").Value = 5026.2
$ws.Range("K134").Value = 9958.363799999999
$ws.Range("L134").Value = 15078.6
$ws.Range("M134").Value = -7423.363799999999
$ws.Range("N134").Value = -20148.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6116.4814
$ws.Range("I31").Value = 2238.5557
$ws.Range("J31").Value = 8055.4443
$ws.Range("K31").Value = 2238.5557
$ws.Range("L31").Value = 8055.4443
$ws.Range("M31").Value = -1943.5557
$ws.Range("N31").Value = -8645.444299999999

$ws.Range("H34").Value = 6116.4814
$ws.Range("I34").Value = 2238.5557
$ws.Range("J34").Value = 8055.4443
$ws.Range("K34").Value = 2238.5557
$ws.Range("L34").Value = 8055.4443
$ws.Range("M34").Value = -2036.5557
$ws.Range("N34").Value = -8459.444299999999

$ws.Range("H58").Value = 1844.2
$ws.Range("I58").Value = 719.619
$ws.Range("J58").Value = 4468.222
$ws.Range("K58").Value = 719.619
$ws.Range("L58").Value = 4468.222
$ws.Range("M58").Value = -516.619
$ws.Range("N58").Value = -4874.222

$ws.Range("H86").Value = 4600.35
$ws.Range("I86").Value = 4951.1665
$ws.Range("J86").Value = 4450
$ws.Range("K86").Value = 4951.1665
$ws.Range("L86").Value = 4450
$ws.Range("M86").Value = -3828.1665
$ws.Range("N86").Value = -6696

$ws.Range("H89").Value = 4600.35
$ws.Range("I89").Value = 4951.1665
$ws.Range("J89").Value = 4450
$ws.Range("K89").Value = 24755.8325
$ws.Range("L89").Value = 22250
$ws.Range("M89").Value = -19139.8325
$ws.Range("N89").Value = -33482

$ws.Range("H105").Value = 1665.7391
$ws.Range("I105").Value = 1506
$ws.Range("J105").Value = 2424.5
$ws.Range("K105").Value = 1506
$ws.Range("L105").Value = 2424.5
$ws.Range("M105").Value = 241
$ws.Range("N105").Value = -5918.5

$ws.Range("H136").Value = 1844.2
$ws.Range("I136").Value = 719.619
$ws.Range("J136").Value = 4468.222
$ws.Range("K136").Value = 2158.857
$ws.Range("L136").Value = 13404.666
$ws.Range("M136").Value = 391.143
$ws.Range("N136").Value = -18504.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1093.7894
$ws.Range("I5").Value = 287.44446
$ws.Range("J5").Value = 1819.5
$ws.Range("K5").Value = 862.33338
$ws.Range("L5").Value = 5458.5
$ws.Range("M5").Value = -750.33338
$ws.Range("N5").Value = -5682.5

$ws.Range("H12").Value = 32.666668
$ws.Range("J12").Value = 26.928572
$ws.Range("L12").Value = 80.78571599999999
$ws.Range("N12").Value = -426.785716

$ws.Range("H75").Value = 5001
$ws.Range("J75").Value = 5001
$ws.Range("L75").Value = 15003
$ws.Range("N75").Value = -16999

$ws.Range("H78").Value = 5001
$ws.Range("J78").Value = 5001
$ws.Range("L78").Value = 45009
$ws.Range("N78").Value = -54993

$ws.Range("H135").Value = 1093.7894
$ws.Range("I135").Value = 287.44446
$ws.Range("J135").Value = 1819.5
$ws.Range("K135").Value = 2587.00014
$ws.Range("L135").Value = 16375.5
$ws.Range("M135").Value = -52.0001400000001
$ws.Range("N135").Value = -21445.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 10620.444
$ws.Range("J94").Value = 10620.444
$ws.Range("L94").Value = 10620.444
$ws.Range("N94").Value = -11972.444

$ws.Range("H97").Value = 1147.9565
$ws.Range("I97").Value = 912.1429000000001
$ws.Range("J97").Value = 1514.7778
$ws.Range("K97").Value = 912.1429000000001
$ws.Range("L97").Value = 1514.7778
$ws.Range("M97").Value = -416.1429000000001
$ws.Range("N97").Value = -2506.7778

$ws.Range("H132").Value = 9208.706
$ws.Range("I132").Value = 9770.134
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 29310.402
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -26780.402
$ws.Range("N132").Value = -20054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4563.4707
$ws.Range("I40").Value = 5789.091
$ws.Range("J40").Value = 2316.5
$ws.Range("K40").Value = 5789.091
$ws.Range("L40").Value = 2316.5
$ws.Range("M40").Value = -5653.091
$ws.Range("N40").Value = -2588.5

$ws.Range("H61").Value = 1293.125
$ws.Range("I61").Value = 997.25
$ws.Range("J61").Value = 1884.875
$ws.Range("K61").Value = 997.25
$ws.Range("L61").Value = 1884.875
$ws.Range("M61").Value = -795.25
$ws.Range("N61").Value = -2288.875

$ws.Range("H113").Value = 1293.125
$ws.Range("I113").Value = 997.25
$ws.Range("J113").Value = 1884.875
$ws.Range("K113").Value = 997.25
$ws.Range("L113").Value = 1884.875
$ws.Range("M113").Value = 1172.75
$ws.Range("N113").Value = -6224.875

$ws.Range("H132").Value = 11118706
$ws.Range("I132").Value = 5647
$ws.Range("J132").Value = 38473930
$ws.Range("K132").Value = 16941
$ws.Range("L132").Value = 115421790
$ws.Range("M132").Value = -14411
$ws.Range("N132").Value = -115426850

$ws.Range("H136").Value = 4361.2554
$ws.Range("I136").Value = 1899.1842
$ws.Range("J136").Value = 14756.667
$ws.Range("K136").Value = 5697.5526
$ws.Range("L136").Value = 44270.001
$ws.Range("M136").Value = -3147.5526
$ws.Range("N136").Value = -49370.001
